$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$np = $s.NotesPage
Write-Host ($np | Get-Member | Out-String)
